# Daily attendance processing - 2025-10-05 21:15:12
# Rotate the "Recorded By" (column G) contributor list so that the last
# recorder listed is moved to the front, for every row that has more than
# one recorder (i.e. a comma-separated list of names/emails).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        if ($trimmed.Count -gt 1) {
            $last = $trimmed[$trimmed.Count - 1]
            $rest = $trimmed[0..($trimmed.Count - 2)]
            $rotated = @($last) + $rest
            $cell.Value2 = [string]::Join(", ", $rotated)
        }
    }
}
